$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) formatting, used to restore
# style on Price cells after forcing a Text number format so that
# numeric-looking strings (e.g. "211.51") are not coerced into floats.
$defaultStyle = $ws.Range("D4").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.658.99"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.598.64"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.51"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0618"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.60"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0839"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.823.26"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.574.01"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.86"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.650.97"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0735"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "207.90"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("E21").Value = "  +5.98%  "
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.32"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.12"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.30"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +2.08%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.281.02"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.621"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -9.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("E39").Value = "  +19.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.836"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  +3.15%  "
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.786"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.92"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.735.46"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.39"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.101"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  +3.33%  "
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.41"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -1.28%  "
